$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Motherboard")
$ws.Range("A1").Value = "TEST"
